# adloori to davuluri completed
# Fill in the "Total Points" (column E) awarded scores for the Generic
# and Customer Class sections to match the rubric (column D), and move
# the active selection to E15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generic section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Customer Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move/restore the view so the active cell/selection is E15
$ws.Range("E15").Select()
